$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename use-case UC.007 -> UC.006 "Iniciar Contrato" (cell C8 of the Product Backlog table)
$ws.Range("C8").Value = "UC.006 Iniciar Contrato"

# Center the header row of the backlog table (row 2: ITEM BP / REQUISITOS / DESCRIÇÃO / ESTIMATIVA)
$ws.Range("B2:E2").HorizontalAlignment = -4108

# Center the "REQUISITOS" column contents (column C, rows 3-8)
$ws.Range("C3:C8").HorizontalAlignment = -4108

# Leave the selection where it was left when the workbook was last saved
$ws.Range("G6").Select()
